# Applies the Sun Apr 23 18:35:05 UTC 2023 cryptos-list refresh:
# updates Price/Volume(1h) figures for every coin row, and swaps the
# Litecoin/ShibaInu rows (18 <-> 19) to reflect the new ranking order.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Many 'Price' cells hold plain numeric-looking text (e.g. '86.83').
# Assigning such a string straight to .Value lets Excel auto-coerce it
# into a real number, which would change the cell's stored type.
# Writing it as a quoted-text formula and then collapsing the formula
# to its value via Copy + PasteSpecial(xlPasteValues) keeps the cell a
# plain text value (matching the source data) without touching styles.
function Set-TextCell($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
}

$updates = @(
    @{ Ref = 'D2'; Value = '27.893.30'; ForceText = $false }
    @{ Ref = 'E2'; Value = '  +0.90%  '; ForceText = $false }
    @{ Ref = 'D3'; Value = '1.886.16'; ForceText = $false }
    @{ Ref = 'E3'; Value = '  +0.34%  '; ForceText = $false }
    @{ Ref = 'E4'; Value = '  +1.75%  '; ForceText = $false }
    @{ Ref = 'D5'; Value = '335.72'; ForceText = $true }
    @{ Ref = 'E5'; Value = '  +1.33%  '; ForceText = $false }
    @{ Ref = 'D6'; Value = '1.018'; ForceText = $true }
    @{ Ref = 'E6'; Value = '  +1.60%  '; ForceText = $false }
    @{ Ref = 'D7'; Value = '0.4651'; ForceText = $true }
    @{ Ref = 'E7'; Value = '  -1.54%  '; ForceText = $false }
    @{ Ref = 'D8'; Value = '0.3907'; ForceText = $true }
    @{ Ref = 'E8'; Value = '  -1.38%  '; ForceText = $false }
    @{ Ref = 'D9'; Value = '47.04'; ForceText = $true }
    @{ Ref = 'E9'; Value = '  -0.87%  '; ForceText = $false }
    @{ Ref = 'D10'; Value = '0.07955'; ForceText = $true }
    @{ Ref = 'E10'; Value = '  -0.76%  '; ForceText = $false }
    @{ Ref = 'D11'; Value = '1.010'; ForceText = $true }
    @{ Ref = 'E11'; Value = '  -1.21%  '; ForceText = $false }
    @{ Ref = 'D12'; Value = '21.52'; ForceText = $true }
    @{ Ref = 'E12'; Value = '  -1.00%  '; ForceText = $false }
    @{ Ref = 'D13'; Value = '1.900.08'; ForceText = $false }
    @{ Ref = 'E13'; Value = '  +1.18%  '; ForceText = $false }
    @{ Ref = 'D14'; Value = '5.931'; ForceText = $true }
    @{ Ref = 'E14'; Value = '  -0.52%  '; ForceText = $false }
    @{ Ref = 'D15'; Value = '7.079'; ForceText = $true }
    @{ Ref = 'E15'; Value = '  -1.14%  '; ForceText = $false }
    @{ Ref = 'E16'; Value = '  +1.71%  '; ForceText = $false }
    @{ Ref = 'D17'; Value = '0.06757'; ForceText = $true }
    @{ Ref = 'E17'; Value = '  +2.12%  '; ForceText = $false }
    @{ Ref = 'B18'; Value = 'ShibaInu'; ForceText = $false }
    @{ Ref = 'C18'; Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; ForceText = $false }
    @{ Ref = 'D18'; Value = '0.00001045'; ForceText = $true }
    @{ Ref = 'E18'; Value = '  +0.58%  '; ForceText = $false }
    @{ Ref = 'B19'; Value = 'Litecoin'; ForceText = $false }
    @{ Ref = 'C19'; Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'; ForceText = $false }
    @{ Ref = 'D19'; Value = '86.83'; ForceText = $true }
    @{ Ref = 'E19'; Value = '  -0.02%  '; ForceText = $false }
    @{ Ref = 'D20'; Value = '17.03'; ForceText = $true }
    @{ Ref = 'E20'; Value = '  -0.83%  '; ForceText = $false }
    @{ Ref = 'E21'; Value = '  +1.67%  '; ForceText = $false }
    @{ Ref = 'D22'; Value = '27.918.49'; ForceText = $false }
    @{ Ref = 'E22'; Value = '  +0.91%  '; ForceText = $false }
    @{ Ref = 'D23'; Value = '5.478'; ForceText = $true }
    @{ Ref = 'E23'; Value = '  -0.37%  '; ForceText = $false }
    @{ Ref = 'D24'; Value = '10.90'; ForceText = $true }
    @{ Ref = 'E24'; Value = '  -0.66%  '; ForceText = $false }
    @{ Ref = 'D25'; Value = '2.351'; ForceText = $true }
    @{ Ref = 'E25'; Value = '  +1.70%  '; ForceText = $false }
    @{ Ref = 'D26'; Value = '2.115.96'; ForceText = $false }
    @{ Ref = 'E26'; Value = '  +0.82%  '; ForceText = $false }
    @{ Ref = 'D27'; Value = '159.08'; ForceText = $true }
    @{ Ref = 'E27'; Value = '  +1.93%  '; ForceText = $false }
    @{ Ref = 'D28'; Value = '19.95'; ForceText = $true }
    @{ Ref = 'E28'; Value = '  -1.34%  '; ForceText = $false }
    @{ Ref = 'D29'; Value = '2.063'; ForceText = $true }
    @{ Ref = 'E29'; Value = '  -1.21%  '; ForceText = $false }
    @{ Ref = 'D30'; Value = '5.405'; ForceText = $true }
    @{ Ref = 'E30'; Value = '  -2.76%  '; ForceText = $false }
    @{ Ref = 'D31'; Value = '121.24'; ForceText = $true }
    @{ Ref = 'E31'; Value = '  -0.69%  '; ForceText = $false }
    @{ Ref = 'D32'; Value = '0.9587'; ForceText = $true }
    @{ Ref = 'E32'; Value = '  -0.62%  '; ForceText = $false }
    @{ Ref = 'D33'; Value = '0.09479'; ForceText = $true }
    @{ Ref = 'E33'; Value = '  -0.55%  '; ForceText = $false }
    @{ Ref = 'D34'; Value = '3.671'; ForceText = $true }
    @{ Ref = 'E34'; Value = '  +1.08%  '; ForceText = $false }
    @{ Ref = 'D35'; Value = '1.367'; ForceText = $true }
    @{ Ref = 'E35'; Value = '  -6.20%  '; ForceText = $false }
    @{ Ref = 'D36'; Value = '5.318'; ForceText = $true }
    @{ Ref = 'E36'; Value = '  +0.37%  '; ForceText = $false }
    @{ Ref = 'D37'; Value = '0.06088'; ForceText = $true }
    @{ Ref = 'E37'; Value = '  -0.36%  '; ForceText = $false }
    @{ Ref = 'D38'; Value = '0.02232'; ForceText = $true }
    @{ Ref = 'E38'; Value = '  -0.74%  '; ForceText = $false }
    @{ Ref = 'D39'; Value = '1.217'; ForceText = $true }
    @{ Ref = 'E39'; Value = '  -0.88%  '; ForceText = $false }
    @{ Ref = 'D40'; Value = '8.049'; ForceText = $true }
    @{ Ref = 'E40'; Value = '  -1.10%  '; ForceText = $false }
    @{ Ref = 'D41'; Value = '0.5925'; ForceText = $true }
    @{ Ref = 'E41'; Value = '  -1.19%  '; ForceText = $false }
    @{ Ref = 'D42'; Value = '0.1878'; ForceText = $true }
    @{ Ref = 'E42'; Value = '  -0.78%  '; ForceText = $false }
    @{ Ref = 'E43'; Value = '  -0.33%  '; ForceText = $false }
    @{ Ref = 'D44'; Value = '1.272'; ForceText = $true }
    @{ Ref = 'E44'; Value = '  +1.12%  '; ForceText = $false }
    @{ Ref = 'D45'; Value = '0.5645'; ForceText = $true }
    @{ Ref = 'E45'; Value = '  -1.04%  '; ForceText = $false }
    @{ Ref = 'D46'; Value = '12.10'; ForceText = $true }
    @{ Ref = 'E46'; Value = '  -1.29%  '; ForceText = $false }
    @{ Ref = 'D47'; Value = '3.395'; ForceText = $true }
    @{ Ref = 'E47'; Value = '  -0.31%  '; ForceText = $false }
    @{ Ref = 'D48'; Value = '1.914'; ForceText = $true }
    @{ Ref = 'E48'; Value = '  -0.81%  '; ForceText = $false }
    @{ Ref = 'D49'; Value = '0.06904'; ForceText = $true }
    @{ Ref = 'E49'; Value = '  +1.24%  '; ForceText = $false }
    @{ Ref = 'D50'; Value = '113.58'; ForceText = $true }
    @{ Ref = 'E50'; Value = '  +2.47%  '; ForceText = $false }
    @{ Ref = 'E51'; Value = '  -0.23%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Ref)
    if ($u.ForceText) {
        Set-TextCell $range $u.Value
    } else {
        $range.Value = $u.Value
    }
}
